$d = $word.ActiveDocument

# Convert an RRGGBB hex string into the BGR-packed integer that the
# Word OM's Font.Color property expects.
function ToWordColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$highlightColor = ToWordColor("2C3E50")

# Bold + color a sub-string of a paragraph identified by its (unique)
# full text. Applying Font formatting to a narrower Range than the run
# automatically splits the surrounding run(s) in the saved OOXML - the
# same effect Word itself produces when you select text and click Bold.
function HighlightSubstring($paragraphText, [string[]]$needles) {
    $full = $d.Content.Duplicate
    $full.Find.ClearFormatting()
    $found = $full.Find.Execute($paragraphText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not locate paragraph text: $paragraphText"
    }
    $paraStart = $full.Start

    $searchFrom = 0
    foreach ($needle in $needles) {
        $idx = $paragraphText.IndexOf($needle, $searchFrom)
        if ($idx -lt 0) {
            throw "Could not locate needle '$needle' inside paragraph text"
        }
        $s = $paraStart + $idx
        $e = $s + $needle.Length
        $sub = $d.Range($s, $e)
        $sub.Font.Bold = $true
        $sub.Font.Color = $highlightColor
        $searchFrom = $idx + $needle.Length
    }
}

HighlightSubstring "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%" @("23%", "64%")

HighlightSubstring "• Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes" @("±4.2%", "±2.1%", "71%", "87%")

HighlightSubstring "• Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis" @("73.5%", "`$4.7M")

HighlightSubstring "• Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion" @("`$2")

HighlightSubstring "• Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%" @("57%")

HighlightSubstring "• Revenue generation: Delivered `$4.9M additional revenue through optimization" @("`$4.9M")

HighlightSubstring "• 23% conversion rate improvement" @("23%")

HighlightSubstring "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations" @("12,847")

Write-Output "Highlighting complete"
